$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4500
$ws.Range("J62").Value = 6233.3335
$ws.Range("L62").Value = 6233.3335
$ws.Range("N62").Value = -7481.3335
$ws.Range("H65").Value = 4500
$ws.Range("J65").Value = 6233.3335
$ws.Range("L65").Value = 31166.6675
$ws.Range("N65").Value = -37406.6675
$ws.Range("H116").Value = 275798.78
$ws.Range("I116").Value = 558141.6
$ws.Range("J116").Value = 8316.105
$ws.Range("K116").Value = 558141.6
$ws.Range("L116").Value = 8316.105
$ws.Range("M116").Value = -554699.6
$ws.Range("N116").Value = -15200.105
$ws.Range("H129").Value = 850.76
$ws.Range("I129").Value = 376
$ws.Range("J129").Value = 909.43823
$ws.Range("K129").Value = 1128
$ws.Range("L129").Value = 2728.31469
$ws.Range("M129").Value = 3872
$ws.Range("N129").Value = -12728.31469
$ws.Range("H141").Value = 6699.3613
$ws.Range("I141").Value = 7758.3447
$ws.Range("J141").Value = 2312.1428
$ws.Range("K141").Value = 23275.0341
$ws.Range("L141").Value = 6936.428400000001
$ws.Range("M141").Value = -18095.0341
$ws.Range("N141").Value = -17296.4284

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1099.174
$ws.Range("I2").Value = 1130.1875
$ws.Range("J2").Value = 1028.2858
$ws.Range("K2").Value = 1130.1875
$ws.Range("L2").Value = 1028.2858
$ws.Range("M2").Value = -1017.1875
$ws.Range("N2").Value = -1254.2858
$ws.Range("H32").Value = 4932.857
$ws.Range("I32").Value = 5054.488
$ws.Range("J32").Value = 4600.4
$ws.Range("K32").Value = 5054.488
$ws.Range("L32").Value = 4600.4
$ws.Range("M32").Value = -4767.488
$ws.Range("N32").Value = -5174.4
$ws.Range("H45").Value = 1123.0588
$ws.Range("I45").Value = 1150.2307
$ws.Range("K45").Value = 1150.2307
$ws.Range("M45").Value = -773.2307000000001
$ws.Range("H74").Value = 4287.1377
$ws.Range("I74").Value = 4805.316
$ws.Range("J74").Value = 3302.6
$ws.Range("K74").Value = 4805.316
$ws.Range("L74").Value = 3302.6
$ws.Range("M74").Value = -3931.316
$ws.Range("N74").Value = -5050.6
$ws.Range("H77").Value = 4287.1377
$ws.Range("I77").Value = 4805.316
$ws.Range("J77").Value = 3302.6
$ws.Range("K77").Value = 24026.58
$ws.Range("L77").Value = 16513
$ws.Range("M77").Value = -19658.58
$ws.Range("N77").Value = -25249
$ws.Range("H116").Value = 1099.174
$ws.Range("I116").Value = 1130.1875
$ws.Range("J116").Value = 1028.2858
$ws.Range("K116").Value = 1130.1875
$ws.Range("L116").Value = 1028.2858
$ws.Range("M116").Value = 1163.8125
$ws.Range("N116").Value = -5616.2858
$ws.Range("H122").Value = 2139.2273
$ws.Range("I122").Value = 1434.125
$ws.Range("K122").Value = 4302.375
$ws.Range("M122").Value = -1852.375

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1099.174
$ws.Range("I3").Value = 1130.1875
$ws.Range("J3").Value = 1028.2858
$ws.Range("K3").Value = 1130.1875
$ws.Range("L3").Value = 1028.2858
$ws.Range("M3").Value = -1016.1875
$ws.Range("N3").Value = -1256.2858

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 35000
$ws.Range("J106").Value = 35000
$ws.Range("L106").Value = 35000
$ws.Range("N106").Value = -37524
$ws.Range("H132").Value = 2630.3076
$ws.Range("I132").Value = 1767.9474
$ws.Range("J132").Value = 4971
$ws.Range("K132").Value = 5303.8422
$ws.Range("L132").Value = 14913
$ws.Range("M132").Value = -2773.8422
$ws.Range("N132").Value = -19973

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1725.1666
$ws.Range("I17").Value = 800
$ws.Range("J17").Value = 2187.75
$ws.Range("K17").Value = 2400
$ws.Range("L17").Value = 6563.25
$ws.Range("M17").Value = -2231
$ws.Range("N17").Value = -6901.25
$ws.Range("H19").Value = 1000
$ws.Range("J19").Value = 1000
$ws.Range("L19").Value = 3000
$ws.Range("N19").Value = -3348
$ws.Range("H102").Value = 1000
$ws.Range("I102").Value = 1000
$ws.Range("K102").Value = 3000
$ws.Range("M102").Value = -566
$ws.Range("H140").Value = 2306.348
$ws.Range("I140").Value = 2539.2632
$ws.Range("K140").Value = 7617.7896
$ws.Range("M140").Value = -2437.7896

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 33068.75
$ws.Range("J88").Value = 33068.75
$ws.Range("L88").Value = 33068.75
$ws.Range("N88").Value = -33970.75
$ws.Range("H91").Value = 33068.75
$ws.Range("J91").Value = 33068.75
$ws.Range("L91").Value = 33068.75
$ws.Range("N91").Value = -36188.75
$ws.Range("H122").Value = 2699.4
$ws.Range("I122").Value = 2027.1111
$ws.Range("K122").Value = 6081.3333
$ws.Range("M122").Value = -3631.3333
$ws.Range("H132").Value = 2493.7812
$ws.Range("I132").Value = 1932.9333
$ws.Range("J132").Value = 2988.647
$ws.Range("K132").Value = 5798.7999
$ws.Range("L132").Value = 8965.940999999999
$ws.Range("M132").Value = -3268.7999
$ws.Range("N132").Value = -14025.941
$ws.Range("H141").Value = 42285.8
$ws.Range("J141").Value = 42285.8
$ws.Range("L141").Value = 42285.8
$ws.Range("N141").Value = -52645.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 28000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 28000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 28000
$ws.Range("N38").Value = -28820
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("H122").Value = 2795.2104
$ws.Range("I122").Value = 1669.5385
$ws.Range("J122").Value = 5234.1665
$ws.Range("K122").Value = 5008.6155
$ws.Range("L122").Value = 15702.4995
$ws.Range("M122").Value = -2558.6155
$ws.Range("N122").Value = -20602.4995
$ws.Range("H136").Value = 3755.111
$ws.Range("I136").Value = 1345.8667
$ws.Range("J136").Value = 6766.6665
$ws.Range("K136").Value = 4037.6001
$ws.Range("L136").Value = 20299.9995
$ws.Range("M136").Value = -1487.6001
$ws.Range("N136").Value = -25399.9995
$ws.Range("M38").ClearContents()
$ws.Range("N103").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4654.154
$ws.Range("I122").Value = 2563
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 7689
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -5239
$ws.Range("N122").Value = -28900
$ws.Range("H126").Value = 2275.8262
$ws.Range("I126").Value = 1595.7693
$ws.Range("J126").Value = 3159.9
$ws.Range("K126").Value = 4787.3079
$ws.Range("L126").Value = 9479.700000000001
$ws.Range("M126").Value = -2317.3079
$ws.Range("N126").Value = -14419.7
